# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets, which carry duplicate copies of the same exhibition rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4698
    $ws.Range("F3").Value = 138
    $ws.Range("F4").Value = 818
}
